# Insert a new data row at row 33 (weekly Fruta/hortaliza update),
# which pushes the existing rows 33-48 down to 34-49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(33).Insert()

$ws.Cells.Item(33, 1).Value = 10
$ws.Cells.Item(33, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value = "La Araucanía"
$ws.Cells.Item(33, 4).Value = 44606
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100107
$ws.Cells.Item(33, 8).Value = "Otros"
$ws.Cells.Item(33, 9).Value = 100107011
$ws.Cells.Item(33, 10).Value = "Tuna"
$ws.Cells.Item(33, 11).Value = "Sin especificar"
$ws.Cells.Item(33, 12).Value = "Segunda"
$ws.Cells.Item(33, 13).Value = 80
$ws.Cells.Item(33, 14).Value = 10000
$ws.Cells.Item(33, 15).Value = 10000
$ws.Cells.Item(33, 16).Value = 10000
$ws.Cells.Item(33, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(33, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(33, 19).Value = 625
$ws.Cells.Item(33, 20).Value = 16
